$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation so numeric-looking strings (e.g. "1.000", "18.80")
# are not auto-converted to numbers by Excel's smart cell-entry parsing.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.858.51"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "1.901.85"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "0.8005"
$ws.Range("E5").Value = "  +5.61%  "
$ws.Range("D6").Value = "240.73"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.3121"
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("D9").Value = "26.06"
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("D10").Value = "0.06867"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").Value = "0.07973"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.913.01"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "0.7328"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").Value = "5.152"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").Value = "92.41"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "29.872.48"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "13.86"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "5.836"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").Value = "244.03"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "0.000007684"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "2.154.06"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "6.863"
$ws.Range("D25").Value = "167.15"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").Value = "9.164"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").Value = "0.1410"
$ws.Range("E27").Value = "  +7.88%  "
$ws.Range("D28").Value = "18.80"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "2.009"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "1.358"
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("D31").Value = "1.514"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "4.268"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("D33").Value = "0.05528"
$ws.Range("E33").Value = "  +3.10%  "
$ws.Range("D34").Value = "4.047"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "1.251"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").Value = "0.7274"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "2.723"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "0.01916"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").Value = "2.785"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").Value = "6.025"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("D42").Value = "71.70"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "0.8339"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").Value = "1.860"
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").Value = "100.28"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "7.526"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").Value = "9.699"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "2.059.08"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "975.76"
$ws.Range("E50").Value = "  +6.77%  "
$ws.Range("D51").Value = "36.01"
$ws.Range("E51").Value = "  -0.69%  "

# Restore original (unset) cell formatting now that text values are committed,
# so the saved style indices match the source workbook.
$ws.Range("D2:E51").ClearFormats()
